$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 24,6
$data[0,0] = 4.919248675181835
$data[0,1] = 4.444518641865266
$data[0,2] = 16.4883950170072
$data[0,3] = 25.98330712881918
$data[0,4] = 32.51817556220112
$data[0,5] = 14.58649485554098
$data[1,0] = 4.749650641995536
$data[1,1] = 4.476679295345932
$data[1,2] = 15.54930813216279
$data[1,3] = 25.8628325103482
$data[1,4] = 32.20013365927797
$data[1,5] = 14.62942355217465
$data[2,0] = 4.644054114636938
$data[2,1] = 4.497253722014158
$data[2,2] = 14.94854815358411
$data[2,3] = 25.80175334974446
$data[2,4] = 32.02427274474193
$data[2,5] = 14.66096056113957
$data[3,0] = 4.600731489586811
$data[3,1] = 4.505847100731531
$data[3,2] = 14.69792960780875
$data[3,3] = 25.78010867863454
$data[3,4] = 31.95755412377848
$data[3,5] = 14.67510329872941
$data[4,0] = 4.5935225675547
$data[4,1] = 4.507286685816489
$data[4,2] = 14.6559729994206
$data[4,3] = 25.77671059238801
$data[4,4] = 31.94677581217788
$data[4,5] = 14.67752934019327
$data[5,0] = 4.643470923946005
$data[5,1] = 4.497368767277898
$data[5,2] = 14.94519132743234
$data[5,3] = 25.80144830157152
$data[5,4] = 32.02335285616903
$data[5,5] = 14.66114608233373
$data[6,0] = 4.861118969828937
$data[6,1] = 4.45543646128034
$data[6,2] = 16.1697581022852
$data[6,3] = 25.93909585414612
$data[6,4] = 32.40452740869576
$data[6,5] = 14.60021533119003
$data[7,0] = 5.273261840325862
$data[7,1] = 4.379728775653839
$data[7,2] = 18.45690302754691
$data[7,3] = 26.31082894146667
$data[7,4] = 33.30251506233174
$data[7,5] = 14.52229387727257
$data[8,0] = 5.563465931670671
$data[8,1] = 4.328018964123547
$data[8,2] = 20.11216349655488
$data[8,3] = 26.64487798257866
$data[8,4] = 34.04822508425956
$data[8,5] = 14.49099857716713
$data[9,0] = 5.692102092442505
$data[9,1] = 4.305331368372868
$data[9,2] = 20.82353206773012
$data[9,3] = 26.80970058468992
$data[9,4] = 34.40464785571038
$data[9,5] = 14.48251740359033
$data[10,0] = 5.740281110487142
$data[10,1] = 4.296859348277849
$data[10,2] = 21.08696447370156
$data[10,3] = 26.87392345287176
$data[10,4] = 34.54195092691207
$data[10,5] = 14.4801422925193
$data[11,0] = 5.729929276039331
$data[11,1] = 4.298678657700345
$data[11,2] = 21.03049351563482
$data[11,3] = 26.86001220946938
$data[11,4] = 34.5122787590954
$data[11,5] = 14.48061647536069
$data[12,0] = 5.696076733980806
$data[12,1] = 4.304631984446089
$data[12,2] = 20.84532379348005
$data[12,3] = 26.81494825781008
$data[12,4] = 34.41589791680013
$data[12,5] = 14.48230518778797
$data[13,0] = 5.675270405735048
$data[13,1] = 4.308294077786764
$data[13,2] = 20.73112856050771
$data[13,3] = 26.7875794348484
$data[13,4] = 34.35716150938089
$data[13,5] = 14.48344876686517
$data[14,0] = 5.554987133248692
$data[14,1] = 4.32951838940854
$data[14,2] = 20.06483961127597
$data[14,3] = 26.63436184087534
$data[14,4] = 34.02526630191414
$data[14,5] = 14.49166943014764
$data[15,0] = 5.480296705338666
$data[15,1] = 4.342752178229222
$data[15,2] = 19.64546010490858
$data[15,3] = 26.54363439276447
$data[15,4] = 33.82596685421764
$data[15,5] = 14.49819315396954
$data[16,0] = 5.437021119347742
$data[16,1] = 4.350442590225636
$data[16,2] = 19.40032463661223
$data[16,3] = 26.49266384808683
$data[16,4] = 33.71296246979541
$data[16,5] = 14.5024866422991
$data[17,0] = 5.422316056848641
$data[17,1] = 4.353059972669789
$data[17,2] = 19.31665167285689
$data[17,3] = 26.47561568478598
$data[17,4] = 33.67498475310814
$data[17,5] = 14.50403301476275
$data[18,0] = 5.488280651570977
$data[18,1] = 4.341335281354494
$data[18,2] = 19.69050914970961
$data[18,3] = 26.55316717601255
$data[18,4] = 33.84701517805356
$data[18,5] = 14.49744261194699
$data[19,0] = 5.706034849372579
$data[19,1] = 4.302880117529673
$data[19,2] = 20.89987374584484
$data[19,3] = 26.8281359178695
$data[19,4] = 34.4441451191644
$data[19,5] = 14.48178640108065
$data[20,0] = 5.845226879203546
$data[20,1] = 4.278442323150961
$data[20,2] = 21.65561435974893
$data[20,3] = 27.0183600753009
$data[20,4] = 34.84793650686856
$data[20,5] = 14.47643373549416
$data[21,0] = 5.771237281426264
$data[21,1] = 4.29142192822723
$data[21,2] = 21.25541952971247
$data[21,3] = 26.9158866979076
$data[21,4] = 34.63123488784295
$data[21,5] = 14.47884129607636
$data[22,0] = 5.484672149520582
$data[22,1] = 4.341975604548146
$data[22,2] = 19.67015502756954
$data[22,3] = 26.54885369543834
$data[22,4] = 33.8374943129646
$data[22,5] = 14.49778024162866
$data[23,0] = 5.163737942608901
$data[23,1] = 4.399518154319042
$data[23,2] = 17.81018603093844
$data[23,3] = 26.19945639577376
$data[23,4] = 33.04402059477714
$data[23,5] = 14.53885838887229

$ws.Range("C2:H25").Value = $data

$kdata = New-Object 'object[,]' 24,1
$kdata[0,0] = 16.19582613023564
$kdata[1,0] = 15.48873386670668
$kdata[2,0] = 15.04004877232157
$kdata[3,0] = 14.8538124505776
$kdata[4,0] = 14.82269136099103
$kdata[5,0] = 15.03755049587162
$kdata[6,0] = 15.95518112642491
$kdata[7,0] = 17.62965461531941
$kdata[8,0] = 18.77277215515562
$kdata[9,0] = 19.27211270608633
$kdata[10,0] = 19.45810938600309
$kdata[11,0] = 19.41819093961609
$kdata[12,0] = 19.28747737025427
$kdata[13,0] = 19.20700552000778
$kdata[14,0] = 18.73971240533232
$kdata[15,0] = 18.44765723767933
$kdata[16,0] = 18.27773519877791
$kdata[17,0] = 18.21987337824691
$kdata[18,0] = 18.47894868816586
$kdata[19,0] = 19.32595591572021
$kdata[20,0] = 19.8614543507497
$kdata[21,0] = 19.57733710042014
$kdata[22,0] = 18.46480809230224
$kdata[23,0] = 17.1912815158128

$ws.Range("K2:K25").Value = $kdata
